# Append the new "11-11-2025" gold-price row (row 56) to Sheet1, matching
# the existing table's layout/formatting (columns A: date, B: description).
#
# NOTE: Assigning a plain ambiguous-date-looking string (e.g. "11-11-2025")
# straight to Range.Value makes Excel's smart text parser re-interpret it
# as a real date (mm-dd-yyyy) and silently replaces the cell's stored
# content with a date serial number + a brand-new number-format style.
# To keep the cell a literal text shared-string (matching every other
# date cell in column A) we instead enter it as a text-literal formula
# (="11-11-2025"), which evaluates to the plain string without being
# reinterpreted, then convert that formula down to a static value via
# Copy/PasteSpecial(values). This preserves the ambient column style
# (border, General number format) instead of minting a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 56

$dateCell = $ws.Range("A" + $lastRow)
$dateCell.Formula = "=""11-11-2025"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

$priceCell = $ws.Range("B" + $lastRow)
$priceCell.Value = "The price of gold in India today is ₹12,628 per gram for 24 karat gold, ₹11,575 per gram for 22 karat gold and ₹9,471 per gram for 18 karat gold (also called 999 gold)."

$excel.CutCopyMode = 0
